# Correções aos ficheiros existentes, colocação da blueprint final
#
# The "Saved" (POST) endpoint in row 14 no longer needs the productID in the
# URL (it travels in the request body instead), so its path shrinks from
# "saved/{userID}/{productID}" to "saved/{userID}" - matching the text
# already used by the "Saved" (GET) row right below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D14").Value = "saved/{userID}"

# Scroll the sheet back so the view starts at the top-left corner (A1)
# instead of the previously scrolled-to A3, and leave the selection on the
# cell that was just edited.
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("D14").Select()
